# Update the grant history table for the FY2012-2016 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptive text blurb (merged A3:D3) from FY 2011-2016 to FY 2012-2016.
$ws.Range("A3").Value = "This table shows the grant awards and award dollars ACF made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the ACF page of this website."

# Update the table caption (merged A7:C7) from FY 2011-2016 to FY 2012-2016.
$ws.Range("A7").Value = "Grant awards and award dollars ACF made for FY 2012-2016."

# Leave A7:C7 selected, matching the saved state after the edit.
$ws.Range("A7:C7").Select()
